# Adds the new "status", "notes", "build" (and, where missing, "gh_issue" /
# "related" / "range_structures") tracking columns to every sheet of the
# model-sets workbook, appending them after the existing header columns.
# This mirrors the "setting up travis full workflow" commit, which extended
# each sheet's schema with build/status bookkeeping columns.

$wb = $excel.ActiveWorkbook

# --- model: name, summary, description, gh_issue, depends_on, +status, +notes, +build
$ws = $wb.Worksheets.Item("model")
$ws.Cells.Item(1, 6).Value = "status"
$ws.Cells.Item(1, 7).Value = "notes"
$ws.Cells.Item(1, 8).Value = "build"
# the existing "current" data row picks up the new status column value and
# leaves notes/build blank, matching the target row (F2="current", G2/H2 empty)
$ws.Cells.Item(2, 6).Value = "current"

# --- packages: name, summary, description, +gh_issue, +status, +notes, +build
$ws = $wb.Worksheets.Item("packages")
$ws.Cells.Item(1, 4).Value = "gh_issue"
$ws.Cells.Item(1, 5).Value = "status"
$ws.Cells.Item(1, 6).Value = "notes"
$ws.Cells.Item(1, 7).Value = "build"

# --- concepts: package, name, summary, description, gh_issue, parents, +related, +status, +notes, +build
$ws = $wb.Worksheets.Item("concepts")
$ws.Cells.Item(1, 7).Value = "related"
$ws.Cells.Item(1, 8).Value = "status"
$ws.Cells.Item(1, 9).Value = "notes"
$ws.Cells.Item(1, 10).Value = "build"

# --- elements: package, name, summary, description, gh_issue, parent, concepts, domains, ranges, +related, +status, +notes, +build
$ws = $wb.Worksheets.Item("elements")
$ws.Cells.Item(1, 10).Value = "related"
$ws.Cells.Item(1, 11).Value = "status"
$ws.Cells.Item(1, 12).Value = "notes"
$ws.Cells.Item(1, 13).Value = "build"

# --- structures: package, name, attribute, element, summary, description, gh_issue, concepts, ranges, +range_structures, +status, +notes, +build
$ws = $wb.Worksheets.Item("structures")
$ws.Cells.Item(1, 10).Value = "range_structures"
$ws.Cells.Item(1, 11).Value = "status"
$ws.Cells.Item(1, 12).Value = "notes"
$ws.Cells.Item(1, 13).Value = "build"
